$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07226
$ws.Range("C2").Value = 1.1002896

$ws.Range("B3").Value = 0.08574
$ws.Range("C3").Value = 1.3029642

$ws.Range("B4").Value = 0.09876
$ws.Range("C4").Value = 1.4990661

$ws.Range("B5").Value = 0.10874
$ws.Range("C5").Value = 1.6953642

$ws.Range("B6").Value = 0.12346
$ws.Range("C6").Value = 1.8901908
